# Sprint 5 product backlog update:
# A handful of backlog items had their Priority (column A) re-scored,
# then the whole backlog (rows 2-19) was re-sorted by Priority ascending
# (stable sort, so items with equal priority keep their relative order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-score the priority column for the affected backlog items ---
$ws.Range("A11").Value = 0
$ws.Range("A12").Value = 0
$ws.Range("A14").Value = 0
$ws.Range("A17").Value = 1
$ws.Range("A18").Value = 1

# --- Re-sort the backlog (A2:D19) by Priority (column A), ascending ---
$sortRange = $ws.Range("A2:D19")
$keyRange = $ws.Range("A2:A19")
$sortRange.Sort($keyRange, 1, $null, $null, 1, $null, 1, 1)

# --- Restore the last active selection ---
$ws.Range("C20").Select()
